# Pipeline finally worked - fix the "cep290_unkown" typo (-> "cep290_unknown")
# on the genotype sheet, then leave the genotype tab active with B2 selected
# (the temperature tab was the previously-active tab and loses that status).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("genotype")

# Correct the misspelled genotype label across the data block (B2:M9).
# "failed" entries (e.g. C10) are untouched.
$ws.Range("B2:M9").Value = "cep290_unknown"

# Bring the genotype sheet to the front and leave B2 as the active cell,
# matching the reviewer's final on-screen state.
[void]$ws.Activate()
[void]$ws.Range("B2").Select()
